$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 55.923077
$ws.Range("I15").Value = 55.923077
$ws.Range("K15").Value = 167.769231
$ws.Range("M15").Value = 1.230769000000009
$ws.Range("H17").Value = 2124.5
$ws.Range("J17").Value = 2332.6667
$ws.Range("L17").Value = 6998.000100000001
$ws.Range("N17").Value = -7334.000100000001
$ws.Range("H38").Value = 66
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null
$ws.Range("H51").Value = 7998.3335
$ws.Range("I51").Value = 7998.3335
$ws.Range("K51").Value = 7998.3335
$ws.Range("M51").Value = -7514.3335
$ws.Range("H70").Value = 15689.8
$ws.Range("J70").Value = 17277.555
$ws.Range("L70").Value = 51832.665
$ws.Range("N70").Value = -52372.665
$ws.Range("H73").Value = 15689.8
$ws.Range("J73").Value = 17277.555
$ws.Range("L73").Value = 51832.665
$ws.Range("N73").Value = -53704.665
$ws.Range("H103").Value = 1226.6666
$ws.Range("I103").Value = 1420
$ws.Range("J103").Value = 840
$ws.Range("K103").Value = 4260
$ws.Range("L103").Value = 2520
$ws.Range("M103").Value = -3674
$ws.Range("N103").Value = -3692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 2668
$ws.Range("I10").Value = 2668
$ws.Range("K10").Value = 2668
$ws.Range("M10").Value = -2498
$ws.Range("H45").Value = 1750
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 3141.1667
$ws.Range("I61").Value = 3141.1667
$ws.Range("K61").Value = 3141.1667
$ws.Range("M61").Value = -2929.1667
$ws.Range("H97").Value = 1734.1428
$ws.Range("I97").Value = 1231.5
$ws.Range("K97").Value = 1231.5
$ws.Range("M97").Value = -735.5
$ws.Range("H122").Value = 1477.4445
$ws.Range("I122").Value = 1399.625
$ws.Range("K122").Value = 4198.875
$ws.Range("M122").Value = -1748.875
$ws.Range("H136").Value = 3141.1667
$ws.Range("I136").Value = 3141.1667
$ws.Range("K136").Value = 9423.500100000001
$ws.Range("M136").Value = -6873.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 455063.9
$ws.Range("I4").Value = 455063.9
$ws.Range("K4").Value = 1365191.7
$ws.Range("M4").Value = -1365079.7
$ws.Range("H14").Value = 1174.5
$ws.Range("I14").Value = 1174.5
$ws.Range("K14").Value = 3523.5
$ws.Range("M14").Value = -3350.5
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null
$ws.Range("H50").Value = 3333.3333
$ws.Range("I50").Value = 500
$ws.Range("J50").Value = 4750
$ws.Range("K50").Value = 1500
$ws.Range("L50").Value = 14250
$ws.Range("M50").Value = -1019
$ws.Range("N50").Value = -15212
$ws.Range("H53").Value = 3333.3333
$ws.Range("I53").Value = 500
$ws.Range("J53").Value = 4750
$ws.Range("K53").Value = 1500
$ws.Range("L53").Value = 14250
$ws.Range("M53").Value = -1019
$ws.Range("N53").Value = -15212
$ws.Range("H82").Value = 2998.5
$ws.Range("I82").Value = 2998.5
$ws.Range("K82").Value = 8995.5
$ws.Range("M82").Value = -8589.5
$ws.Range("H85").Value = 2998.5
$ws.Range("I85").Value = 2998.5
$ws.Range("K85").Value = 8995.5
$ws.Range("M85").Value = -7591.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1827.091
$ws.Range("I102").Value = 1932.8
$ws.Range("J102").Value = 770
$ws.Range("K102").Value = 1932.8
$ws.Range("L102").Value = 770
$ws.Range("M102").Value = -310.8
$ws.Range("N102").Value = -4014
$ws.Range("H132").Value = 6999.6
$ws.Range("I132").Value = 6666.222
$ws.Range("K132").Value = 19998.666
$ws.Range("M132").Value = -17468.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6928
$ws.Range("I7").Value = 7399.4
$ws.Range("J7").Value = 5749.5
$ws.Range("K7").Value = 7399.4
$ws.Range("L7").Value = 5749.5
$ws.Range("M7").Value = -7287.4
$ws.Range("N7").Value = -5973.5
$ws.Range("H22").Value = 1199.84
$ws.Range("I22").Value = 1216.5238
$ws.Range("J22").Value = 1112.25
$ws.Range("K22").Value = 1216.5238
$ws.Range("L22").Value = 1112.25
$ws.Range("M22").Value = -921.5237999999999
$ws.Range("N22").Value = -1702.25
$ws.Range("H27").Value = 1199.84
$ws.Range("I27").Value = 1216.5238
$ws.Range("J27").Value = 1112.25
$ws.Range("K27").Value = 1216.5238
$ws.Range("L27").Value = 1112.25
$ws.Range("M27").Value = -1109.5238
$ws.Range("N27").Value = -1326.25
$ws.Range("H46").Value = 4225.4424
$ws.Range("I46").Value = 1745.4445
$ws.Range("K46").Value = 1745.4445
$ws.Range("M46").Value = -1557.4445
$ws.Range("H55").Value = 3993.625
$ws.Range("I55").Value = 2790.4
$ws.Range("K55").Value = 2790.4
$ws.Range("M55").Value = -2617.4
$ws.Range("H82").Value = 1657.0526
$ws.Range("I82").Value = 1705.625
$ws.Range("J82").Value = 1621.7273
$ws.Range("K82").Value = 1705.625
$ws.Range("L82").Value = 1621.7273
$ws.Range("M82").Value = -1344.625
$ws.Range("N82").Value = -2343.7273
$ws.Range("H85").Value = 1657.0526
$ws.Range("I85").Value = 1705.625
$ws.Range("J85").Value = 1621.7273
$ws.Range("K85").Value = 1705.625
$ws.Range("L85").Value = 1621.7273
$ws.Range("M85").Value = -457.625
$ws.Range("N85").Value = -4117.7273
$ws.Range("H93").Value = 500
$ws.Range("I93").Value = 500
$ws.Range("K93").Value = 500
$ws.Range("M93").Value = 748
$ws.Range("H122").Value = 4783.375
$ws.Range("I122").Value = 4783.375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14350.125
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11900.125
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 6928
$ws.Range("I126").Value = 7399.4
$ws.Range("J126").Value = 5749.5
$ws.Range("K126").Value = 22198.2
$ws.Range("L126").Value = 17248.5
$ws.Range("M126").Value = -19728.2
$ws.Range("N126").Value = -22188.5
$ws.Range("H136").Value = 3749.923
$ws.Range("I136").Value = 3776.111
$ws.Range("K136").Value = 11328.333
$ws.Range("M136").Value = -8778.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H81").Value = 22722.572
$ws.Range("I81").Value = 25676.334
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 51352.668
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -50291.668
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 22722.572
$ws.Range("I84").Value = 25676.334
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 256763.34
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -251459.34
$ws.Range("N84").Value = -60608
